# Weekly update: insert a new price record for "Ají" (Inferno / Primera)
# at row 60 of the "Feria Lagunitas de Puerto Montt" sheet. Inserting the
# row shifts all existing records from row 60 onward down by one (old row
# 60 -> new row 61, ..., old row 145 -> new row 146), matching the target
# diff, and the sheet's used range grows from A1:R145 to A1:R146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 60-145 down to 61-146.
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with this week's record.
$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 44467
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = 100112021
$ws.Range("G60").Value = "Ají"
$ws.Range("H60").Value = "Inferno"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 120
$ws.Range("K60").Value = 48000
$ws.Range("L60").Value = 48000
$ws.Range("M60").Value = 48000
$ws.Range("N60").Value = "$/caja 12 kilos"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 4000
$ws.Range("Q60").Value = 12
$ws.Range("R60").Value = "Hortaliza"
